$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-32 and add new rows 33-37 with shifted/updated data
$ws.Range("A2").Value = 44652
$ws.Range("B2").Value = 6.704599999999999
$ws.Range("C2").Value = 4.833069801330566
$ws.Range("A3").Value = 44682
$ws.Range("B3").Value = 8.163476190476191
$ws.Range("C3").Value = 4.837417125701904
$ws.Range("A4").Value = 44713
$ws.Range("B4").Value = 7.597904761904761
$ws.Range("C4").Value = 4.856917858123779
$ws.Range("A5").Value = 44743
$ws.Range("B5").Value = 7.186949999999999
$ws.Range("C5").Value = 4.413410663604736
$ws.Range("A6").Value = 44774
$ws.Range("B6").Value = 8.779478260869567
$ws.Range("C6").Value = 4.411562919616699
$ws.Range("A7").Value = 44805
$ws.Range("B7").Value = 7.757523809523809
$ws.Range("C7").Value = 4.813495635986328
$ws.Range("A8").Value = 44835
$ws.Range("B8").Value = 6.084904761904762
$ws.Range("C8").Value = 4.734267711639404
$ws.Range("A9").Value = 44866
$ws.Range("B9").Value = 6.429761904761905
$ws.Range("C9").Value = 4.617533206939697
$ws.Range("A10").Value = 44896
$ws.Range("B10").Value = 5.768047619047617
$ws.Range("C10").Value = 3.86989426612854
$ws.Range("A11").Value = 44927
$ws.Range("B11").Value = 3.4228
$ws.Range("C11").Value = 4.314633846282959
$ws.Range("A12").Value = 44958
$ws.Range("B12").Value = 2.437473684210526
$ws.Range("C12").Value = 3.448765277862549
$ws.Range("A13").Value = 44986
$ws.Range("B13").Value = 2.407782608695652
$ws.Range("C13").Value = 2.255731582641602
$ws.Range("A14").Value = 45017
$ws.Range("B14").Value = 2.197263157894737
$ws.Range("C14").Value = 2.091178894042969
$ws.Range("A15").Value = 45047
$ws.Range("B15").Value = 2.299318181818181
$ws.Range("C15").Value = 2.199604749679565
$ws.Range("A16").Value = 45078
$ws.Range("B16").Value = 2.474619047619047
$ws.Range("C16").Value = 2.269962072372437
$ws.Range("A17").Value = 45108
$ws.Range("B17").Value = 2.63655
$ws.Range("C17").Value = 2.537910223007202
$ws.Range("A18").Value = 45139
$ws.Range("B18").Value = 2.645130434782609
$ws.Range("C18").Value = 2.475393533706665
$ws.Range("A19").Value = 45170
$ws.Range("B19").Value = 2.69565
$ws.Range("C19").Value = 2.671053171157837
$ws.Range("A20").Value = 45200
$ws.Range("B20").Value = 3.149181818181818
$ws.Range("C20").Value = 2.671859264373779
$ws.Range("A21").Value = 45231
$ws.Range("B21").Value = 3.055523809523809
$ws.Range("C21").Value = 3.113406658172607
$ws.Range("A22").Value = 45261
$ws.Range("B22").Value = 2.53885
$ws.Range("C22").Value = 3.067882537841797
$ws.Range("A23").Value = 45292
$ws.Range("B23").Value = 2.715
$ws.Range("C23").Value = 2.373803615570068
$ws.Range("A24").Value = 45323
$ws.Range("B24").Value = 1.7955
$ws.Range("C24").Value = 2.605878591537476
$ws.Range("A25").Value = 45352
$ws.Range("B25").Value = 1.7473
$ws.Range("C25").Value = 1.92271876335144
$ws.Range("A26").Value = 45383
$ws.Range("B26").Value = 1.791227272727273
$ws.Range("C26").Value = 1.836468577384949
$ws.Range("A27").Value = 45413
$ws.Range("B27").Value = 2.418
$ws.Range("C27").Value = 1.997172355651855
$ws.Range("A28").Value = 45444
$ws.Range("B28").Value = 2.809578947368421
$ws.Range("C28").Value = 2.319653272628784
$ws.Range("A29").Value = 45474
$ws.Range("B29").Value = 2.208681818181818
$ws.Range("C29").Value = 2.780210256576538
$ws.Range("A30").Value = 45505
$ws.Range("B30").Value = 2.086782608695652
$ws.Range("C30").Value = 2.269118547439575
$ws.Range("A31").Value = 45536
$ws.Range("B31").Value = 2.409250000000001
$ws.Range("C31").Value = 2.346853017807007
$ws.Range("A32").Value = 45566
$ws.Range("B32").Value = 2.576956521739131
$ws.Range("C32").Value = 2.295066356658936
$ws.Range("A33").Value = 45597
$ws.Range("B33").Value = 2.982
$ws.Range("C33").Value = 2.532831192016602
$ws.Range("A34").Value = 45627
$ws.Range("B34").Value = 3.406619047619048
$ws.Range("C34").Value = 2.603587865829468
$ws.Range("A35").Value = 45658
$ws.Range("B35").Value = 3.721380952380952
$ws.Range("C35").Value = 2.969963312149048
$ws.Range("A36").Value = 45689
$ws.Range("B36").Value = 3.740947368421053
$ws.Range("C36").Value = 3.936691284179688
$ws.Range("A37").Value = 45717
$ws.Range("B37").Value = 4.137476190476191
$ws.Range("C37").Value = 4.301484107971191

# Ensure the newly added date cells (rows 33-37) use the same date number format as the existing column A cells
$ws.Range("A33:A37").NumberFormat = $ws.Range("A2").NumberFormat

